$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add three new worksheets, in order, at the end of the workbook:
#   CypherOutput_Message, StatOuput, StatOuput_Message
# ------------------------------------------------------------------
$wsCypherMsg = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsCypherMsg.Name = "CypherOutput_Message"

$wsStat = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsStat.Name = "StatOuput"

$wsStatMsg = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsStatMsg.Name = "StatOuput_Message"

# ------------------------------------------------------------------
# Helper: write a 10-row "connection message" block (Neo4j URL,
# user, password, cypher text, output path) starting at a given row
# of a given worksheet. Mirrors the existing "Message" sheet layout.
# ------------------------------------------------------------------
function Write-MessageBlock($sheet, $startRow, $cypherText) {
    $sheet.Range("A" + $startRow).Value = 'Neo4j_URL:'
    $sheet.Range("A" + ($startRow + 1)).Value = 'bolt://ncidb-q325-c.nci.nih.gov:7687'
    $sheet.Range("A" + ($startRow + 2)).Value = 'User_name:'
    $sheet.Range("A" + ($startRow + 3)).Value = 'neo4j'
    $sheet.Range("A" + ($startRow + 4)).Value = 'PWD:'
    $sheet.Range("A" + ($startRow + 5)).Value = 'icdcDBneo4j0'
    $sheet.Range("A" + ($startRow + 6)).Value = 'Cypher:'
    $sheet.Range("A" + ($startRow + 7)).Value = $cypherText
    $sheet.Range("A" + ($startRow + 8)).Value = 'Output:'
    $sheet.Range("A" + ($startRow + 9)).Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC01_Trials_Filter_TrialCode-NCIMatch_Neo4jData.xlsx'
}

# CypherOutput_Message: identical 10-row block to the "Message" sheet
# (reuses the same shared strings -- no new strings introduced here)
Write-MessageBlock $wsCypherMsg 1 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE t.clinical_trial_designation IN [''NCI-MATCH''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

# ------------------------------------------------------------------
# StatOuput: header row + one data row of stat-query results.
# Written before StatOuput_Message so new shared strings are
# interned in the same order as the original workbook.
# ------------------------------------------------------------------
$wsStat.Range("A1").Value = "number_of_files"
$wsStat.Range("B1").Value = "number_of_cases"
$wsStat.Range("C1").Value = "number_of_trial"
$wsStat.Range("A2").NumberFormat = "@"
$wsStat.Range("A2").Value = "350"
$wsStat.Range("B2").NumberFormat = "@"
$wsStat.Range("B2").Value = "70"
$wsStat.Range("C2").NumberFormat = "@"
$wsStat.Range("C2").Value = "1"

# ------------------------------------------------------------------
# StatOuput_Message: the same 10-row block twice -- once with the
# original trial-filter cypher, once with the new stat-count cypher
# ------------------------------------------------------------------
Write-MessageBlock $wsStatMsg 1 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE t.clinical_trial_designation IN [''NCI-MATCH''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'
Write-MessageBlock $wsStatMsg 11 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE t.clinical_trial_designation IN [''NCI-MATCH''] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial'

# ------------------------------------------------------------------
# Keep the original first sheet ("CypherOutput") the active/selected
# sheet, same as before the edit -- adding sheets at the end otherwise
# leaves the last-added sheet active.
# ------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()

